$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that are removed/moved in rows 2-7
$ws.Range("G2").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("Y4").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("J7").ClearContents()

# Row 2
$ws.Range("A2").Value = '16 - 30 Days'
$ws.Range("B2").Value = 'Losectil'
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 'Losectil V 20mg Capsule'
$ws.Range("I2").Value = 1

# Row 3
$ws.Range("A3").Value = '16 - 30 Days'
$ws.Range("B3").Value = 'Losectil'
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 'Losectil MUPS 20mg Tablet'
$ws.Range("Y3").Value = 2

# Row 4
$ws.Range("A4").Value = '31 - 60 Days'
$ws.Range("B4").Value = 'Losectil'
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 'Losectil 40mg IV Injection'
$ws.Range("AA4").Value = 8

# Row 5
$ws.Range("A5").Value = '31 - 60 Days'
$ws.Range("B5").Value = 'Topiclo'
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 'Topiclo S 6% 10gm Ointment'
$ws.Range("Q5").Value = 1

# Row 6
$ws.Range("A6").Value = '61 - 90 Days'
$ws.Range("B6").Value = 'Esoral'
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 'Esoral 40mg IV Injection'
$ws.Range("T6").Value = 14

# Row 7
$ws.Range("A7").Value = '61 - 90 Days'
$ws.Range("B7").Value = 'Hairgrow'
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 'Hairgrow 5% Topical Solution 60ml'

# Row 8
$ws.Range("A8").Value = '61 - 90 Days'
$ws.Range("B8").Value = 'Hairgrow'
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = 'Hairgrow 2% Topical Solution 60ml'
$ws.Range("X8").Value = 8

# Row 9
$ws.Range("A9").Value = '61 - 90 Days'
$ws.Range("B9").Value = 'Losectil'
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 'Losectil V 40mg Capsule'

# Row 10
$ws.Range("A10").Value = '61 - 90 Days'
$ws.Range("B10").Value = 'H-QUIN'
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = 'H-QUIN 200mg FC Tablet 20''s'

# Row 11
$ws.Range("A11").Value = '91 - 180 Days'
$ws.Range("B11").Value = 'Esoral'
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = 'Esoral 20mg Capsule 60''s'
$ws.Range("Q11").Value = 1

# Row 12
$ws.Range("A12").Value = '91 - 180 Days'
$ws.Range("B12").Value = 'Orogurd'
$ws.Range("C12").Value = 11
$ws.Range("D12").Value = 'Orogurd 15gm Oral Gel'

# Row 13
$ws.Range("A13").Value = '91 - 180 Days'
$ws.Range("B13").Value = 'Esoral'
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 'Esoral 40mg IV Injection'
$ws.Range("M13").Value = 3
$ws.Range("N13").Value = 27
$ws.Range("Q13").Value = 1
$ws.Range("W13").Value = 22

# Row 14
$ws.Range("A14").Value = '91 - 180 Days'
$ws.Range("B14").Value = 'Remivir'
$ws.Range("C14").Value = 13
$ws.Range("D14").Value = 'Remivir 100mg Lyophilized IV Injection 1''s'
$ws.Range("Q14").Value = 19
$ws.Range("W14").Value = 6

# Row 15
$ws.Range("A15").Value = '91 - 180 Days'
$ws.Range("B15").Value = 'Losectil'
$ws.Range("C15").Value = 14
$ws.Range("D15").Value = 'Losectil V 40mg Capsule'
$ws.Range("AE15").Value = 35

# Row 16
$ws.Range("A16").Value = '91 - 180 Days'
$ws.Range("B16").Value = 'Facid'
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 'Facid 10gm BT Cream'

# Row 17
$ws.Range("A17").Value = '91 - 180 Days'
$ws.Range("B17").Value = 'Mycofin'
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 'Mycofin 5gm Cream'
$ws.Range("J17").Value = 19

# Row 18
$ws.Range("A18").Value = '91 - 180 Days'
$ws.Range("B18").Value = 'Losectil'
$ws.Range("C18").Value = 17
$ws.Range("D18").Value = 'Losectil 40mg Powder for Oral Suspension - 30''s'
$ws.Range("J18").Value = 1

# Row 19
$ws.Range("A19").Value = '91 - 180 Days'
$ws.Range("B19").Value = 'Losectil'
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 'Losectil 40mg IV Injection'
$ws.Range("F19").Value = 181
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = 136

# Row 20
$ws.Range("A20").Value = '91 - 180 Days'
$ws.Range("B20").Value = 'Losectil'
$ws.Range("C20").Value = 19
$ws.Range("D20").Value = 'Losectil 20mg PFS 50''s'
$ws.Range("J20").Value = 1

# Row 21
$ws.Range("A21").Value = '91 - 180 Days'
$ws.Range("B21").Value = 'Panoral'
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 'Panoral 40mg Tablet'
$ws.Range("G21").Value = 18
$ws.Range("I21").Value = 16
$ws.Range("X21").Value = 25
$ws.Range("AE21").Value = 46

# Row 22
$ws.Range("A22").Value = '181 - 210 Days'
$ws.Range("B22").Value = 'Remivir'
$ws.Range("C22").Value = 21
$ws.Range("D22").Value = 'Remivir 100mg Lyophilized IV Injection 1''s'
$ws.Range("I22").Value = 39
$ws.Range("N22").Value = 24
$ws.Range("P22").Value = 24
$ws.Range("U22").Value = 2
$ws.Range("V22").Value = 23
$ws.Range("Z22").Value = 116
$ws.Range("AB22").Value = 40
$ws.Range("AD22").Value = 27
$ws.Range("AE22").Value = 23

# Row 23
$ws.Range("A23").Value = '181 - 210 Days'
$ws.Range("B23").Value = 'Licnil'
$ws.Range("C23").Value = 22
$ws.Range("D23").Value = 'Licnil 117gm Lotion'
$ws.Range("S23").Value = 19

# Row 24
$ws.Range("A24").Value = '181 - 210 Days'
$ws.Range("B24").Value = 'Flucoder'
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 'Flucoder 200mg Capsule - 10''s'
$ws.Range("V24").Value = 21

Write-Output "Update complete"